$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 67, pushing the existing rows 67-116 down to 69-118.
$ws.Rows.Item(67).Resize(2).Insert()

# New row 67 data
$ws.Cells.Item(67, 1).Value = 2
$ws.Cells.Item(67, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(67, 3).Value = "Coquimbo"
$ws.Cells.Item(67, 4).Value = 44574
$ws.Cells.Item(67, 5).Value = 4
$ws.Cells.Item(67, 6).Value = 100112043
$ws.Cells.Item(67, 7).Value = "Pepino ensalada"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 400
$ws.Cells.Item(67, 11).Value = 8500
$ws.Cells.Item(67, 12).Value = 9000
$ws.Cells.Item(67, 13).Value = 8750
$ws.Cells.Item(67, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(67, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(67, 16).Value = 125
$ws.Cells.Item(67, 17).Value = 70
$ws.Cells.Item(67, 18).Value = "Hortaliza"

# New row 68 data
$ws.Cells.Item(68, 1).Value = 2
$ws.Cells.Item(68, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(68, 3).Value = "Coquimbo"
$ws.Cells.Item(68, 4).Value = 44574
$ws.Cells.Item(68, 5).Value = 4
$ws.Cells.Item(68, 6).Value = 100112043
$ws.Cells.Item(68, 7).Value = "Pepino ensalada"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Segunda"
$ws.Cells.Item(68, 10).Value = 240
$ws.Cells.Item(68, 11).Value = 6500
$ws.Cells.Item(68, 12).Value = 7000
$ws.Cells.Item(68, 13).Value = 6750
$ws.Cells.Item(68, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(68, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(68, 16).Value = 68
$ws.Cells.Item(68, 17).Value = 100
$ws.Cells.Item(68, 18).Value = "Hortaliza"
